$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 44
$ws.Range("H44").Value = 0.37
$ws.Range("J44").Value = 18.246575342465754
$ws.Range("K44").Value = 118.83780821917806
$ws.Range("L44").Value = 0.38198888888888854

# Row 45
$ws.Range("H45").Value = 0.37
$ws.Range("J45").Value = 27.36986301369863
$ws.Range("K45").Value = 128.00671232876712
$ws.Range("L45").Value = 0.37860925925925926

# Row 46
$ws.Range("H46").Value = 0.37
$ws.Range("J46").Value = 36.493150684931507
$ws.Range("K46").Value = 137.17561643835614
$ws.Range("L46").Value = 0.37691944444444425

# Row 47
$ws.Range("H47").Value = 0.27
$ws.Range("J47").Value = 13.315068493150687
$ws.Range("K47").Value = 118.98082191780823
$ws.Range("L47").Value = 0.38488888888888911

# Row 48
$ws.Range("H48").Value = 0.3
$ws.Range("J48").Value = 22.191780821917806
$ws.Range("K48").Value = 128.30136986301372
$ws.Range("L48").Value = 0.38259259259259298

# Row 49
$ws.Range("H49").Value = 0.31
$ws.Range("J49").Value = 30.575342465753426
$ws.Range("K49").Value = 137.10410958904109
$ws.Range("L49").Value = 0.37619444444444444

# Row 50
$ws.Range("H50").Value = 0.15
$ws.Range("J50").Value = 7.397260273972603
$ws.Range("K50").Value = 118.13698630136989
$ws.Range("L50").Value = 0.36777777777777826

# Row 51
$ws.Range("H51").Value = 0.21
$ws.Range("J51").Value = 15.534246575342465
$ws.Range("K51").Value = 127.08767123287672
$ws.Range("L51").Value = 0.36618518518518528

# Row 52
$ws.Range("H52").Value = 0.24
$ws.Range("J52").Value = 23.671232876712327
$ws.Range("K52").Value = 136.03835616438357
$ws.Range("L52").Value = 0.3653888888888891

# Row 112
$ws.Range("H112").Value = 0.33
$ws.Range("J112").Value = 5.4246575342465757
$ws.Range("K112").Value = 107.5331506849315
$ws.Range("L112").Value = 0.45826666666666627

# Row 113
$ws.Range("H113").Value = 0.33
$ws.Range("J113").Value = 8.1369863013698627
$ws.Range("K113").Value = 110.29972602739726
$ws.Range("L113").Value = 0.41771111111111114

# Row 114
$ws.Range("H114").Value = 0.33
$ws.Range("J114").Value = 16.273972602739725
$ws.Range("K114").Value = 118.59945205479453
$ws.Range("L114").Value = 0.37715555555555552

# Row 115
$ws.Range("H115").Value = 0.33
$ws.Range("J115").Value = 24.410958904109588
$ws.Range("K115").Value = 126.89917808219177
$ws.Range("L115").Value = 0.36363703703703698

# Row 116
$ws.Range("H116").Value = 0.33
$ws.Range("J116").Value = 32.547945205479451
$ws.Range("K116").Value = 135.19890410958905
$ws.Range("L116").Value = 0.3568777777777779
